$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A widens to fit the new date/time "Bag" timestamp values.
$ws.Columns.Item(1).ColumnWidth = 14

# Apply a date/time number format to column A (style gets reused for A1 header + A2 data).
$ws.Range("A1:A2").NumberFormat = "m/d/yy h:mm"

# New row of sentiment-analysis data (row 2), with a "Bag" method tag in column N.
$ws.Range("A2").Value = 42605.648634259262

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "Bag"
